# SimpleTherapy_Details.xlsx - OTP_Data sheet updates:
#  - top-align the header/value rows (drop wrap text on row 1)
#  - update DOB, Address1 text
#  - add 3 new columns: FirstName_Updated, LastName_Updated, City_Updated
#  - set portrait page orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OTP_Data")
$ws.Activate()

# --- Re-align existing header/data cells to top (removes the old wrap-only look) ---
$ws.Range("A1:K2").VerticalAlignment = -4160   # xlTop

# --- Update DOB (F2) ---
$ws.Range("F2").Value = 34700

# --- Update Address1 (G2) ---
$ws.Range("G2").Value = "1950 San Benito Dr"

# --- New columns: headers ---
$ws.Range("L1").Value = "FirstName_Updated"
$ws.Range("M1").Value = "LastName_Updated"
$ws.Range("N1").Value = "City_Updated"

# --- New columns: values ---
$ws.Range("L2").Value = "Sagar1"
$ws.Range("M2").Value = "Pangale1"
$ws.Range("N2").Value = "Fremont1"

# match header style (top aligned, bold font w/ border) used by the other header cells
$ws.Range("L1:N1").VerticalAlignment = -4160
$ws.Range("L2:N2").VerticalAlignment = -4160

# --- column widths for the new columns ---
$ws.Columns("L:M").ColumnWidth = 17.27
$ws.Columns("N:N").ColumnWidth = 13.27

# --- view tweaks ---
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("K9").Select()

# --- page setup ---
$ws.PageSetup.Orientation = 1   # xlPortrait

Write-Output "done"
